$wb = $excel.ActiveWorkbook

$sheetNames = @("Q2_20_21 all data", "Q1_20_21 all data")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A3").Value = "DfT Group"

    for ($r = 4; $r -le 21; $r++) {
        $ws.Cells.Item($r, 1).Value = "Rail"
    }
}
